$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "28.609.17"
Set-TextValue 2 5 "  +2.03%  "

# Row 3
Set-TextValue 3 4 "1.908.82"
Set-TextValue 3 5 "  +2.43%  "

# Row 4
Set-TextValue 4 4 "1.014"
Set-TextValue 4 5 "  +1.05%  "

# Row 5
Set-TextValue 5 4 "319.23"
Set-TextValue 5 5 "  +2.26%  "

# Row 6
Set-TextValue 6 4 "1.020"
Set-TextValue 6 5 "  +1.70%  "

# Row 7
Set-TextValue 7 4 "0.5210"
Set-TextValue 7 5 "  +2.21%  "

# Row 8
Set-TextValue 8 4 "0.3962"
Set-TextValue 8 5 "  +3.37%  "

# Row 9
Set-TextValue 9 4 "0.08502"
Set-TextValue 9 5 "  +2.36%  "

# Row 10
Set-TextValue 10 4 "1.138"
Set-TextValue 10 5 "  +2.04%  "

# Row 11
Set-TextValue 11 4 "42.16"
Set-TextValue 11 5 "  +1.56%  "

# Row 12
Set-TextValue 12 4 "6.344"
Set-TextValue 12 5 "  +1.90%  "

# Row 13
Set-TextValue 13 4 "20.91"
Set-TextValue 13 5 "  +1.45%  "

# Row 14
Set-TextValue 14 4 "1.876.71"
Set-TextValue 14 5 "  +0.97%  "

# Row 15
Set-TextValue 15 4 "7.389"
Set-TextValue 15 5 "  +2.37%  "

# Row 16
Set-TextValue 16 4 "1.035"
Set-TextValue 16 5 "  +3.17%  "

# Row 17
Set-TextValue 17 4 "0.00001120"
Set-TextValue 17 5 "  +1.98%  "

# Row 18
Set-TextValue 18 4 "92.07"
Set-TextValue 18 5 "  +1.27%  "

# Row 19
Set-TextValue 19 4 "0.06799"
Set-TextValue 19 5 "  +2.48%  "

# Row 20
Set-TextValue 20 4 "18.06"
Set-TextValue 20 5 "  +1.88%  "

# Row 21
Set-TextValue 21 4 "1.017"
Set-TextValue 21 5 "  +1.50%  "

# Row 22
Set-TextValue 22 4 "6.120"
Set-TextValue 22 5 "  +1.30%  "

# Row 23
Set-TextValue 23 4 "28.466.36"
Set-TextValue 23 5 "  +1.44%  "

# Row 24
Set-TextValue 24 4 "11.30"
Set-TextValue 24 5 "  +1.82%  "

# Row 25
Set-TextValue 25 4 "2.291"
Set-TextValue 25 5 "  +2.55%  "

# Row 26
Set-TextValue 26 2 "LidoDAOToken"
Set-TextValue 26 3 "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue 26 4 "2.513"
Set-TextValue 26 5 "  -1.35%  "

# Row 27
Set-TextValue 27 2 "Monero"
Set-TextValue 27 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue 27 4 "160.73"
Set-TextValue 27 5 "  +1.58%  "

# Row 28
Set-TextValue 28 4 "21.01"
Set-TextValue 28 5 "  +2.25%  "

# Row 29
Set-TextValue 29 4 "127.55"
Set-TextValue 29 5 "  +2.22%  "

# Row 30
Set-TextValue 30 4 "0.1059"
Set-TextValue 30 5 "  +0.35%  "

# Row 31
Set-TextValue 31 4 "1.054"
Set-TextValue 31 5 "  +1.20%  "

# Row 32
Set-TextValue 32 4 "5.912"
Set-TextValue 32 5 "  +1.07%  "

# Row 33
Set-TextValue 33 4 "3.659"
Set-TextValue 33 5 "  +1.81%  "

# Row 34
Set-TextValue 34 4 "9.809"
Set-TextValue 34 5 "  +3.40%  "

# Row 35
Set-TextValue 35 4 "0.02474"
Set-TextValue 35 5 "  +2.17%  "

# Row 36
Set-TextValue 36 4 "0.06667"
Set-TextValue 36 5 "  +2.07%  "

# Row 37
Set-TextValue 37 4 "0.2231"
Set-TextValue 37 5 "  +2.59%  "

# Row 38
Set-TextValue 38 4 "1.219"
Set-TextValue 38 5 "  +0.93%  "

# Row 39
Set-TextValue 39 4 "0.6518"
Set-TextValue 39 5 "  +0.77%  "

# Row 40
Set-TextValue 40 4 "1.251"
Set-TextValue 40 5 "  +1.90%  "

# Row 41
Set-TextValue 41 4 "5.042"
Set-TextValue 41 5 "  +1.66%  "

# Row 42
Set-TextValue 42 4 "11.46"
Set-TextValue 42 5 "  +2.14%  "

# Row 43
Set-TextValue 43 4 "0.6204"
Set-TextValue 43 5 "  +1.62%  "

# Row 44
Set-TextValue 44 4 "13.25"
Set-TextValue 44 5 "  +0.98%  "

# Row 45
Set-TextValue 45 4 "1.296"
Set-TextValue 45 5 "  +0.86%  "

# Row 46
Set-TextValue 46 4 "3.724"
Set-TextValue 46 5 "  +1.45%  "

# Row 47
Set-TextValue 47 4 "2.047"
Set-TextValue 47 5 "  +1.46%  "

# Row 48
Set-TextValue 48 4 "1.254"
Set-TextValue 48 5 "  +3.82%  "

# Row 49
Set-TextValue 49 4 "122.24"
Set-TextValue 49 5 "  +1.55%  "

# Row 50
Set-TextValue 50 4 "0.06981"
Set-TextValue 50 5 "  +1.80%  "

# Row 51
Set-TextValue 51 4 "78.64"
Set-TextValue 51 5 "  +0.46%  "
